# Apply the "fix: fixed formatting when scrapping floating point numbers" edit.
#
# The "Importe" column (H) and three "Razon social" entries were scraped with
# Spanish/Argentine number formatting (thousands="." decimal=",") or stray
# commas acting as separators between co-owner names. Both are stored as plain
# TEXT cells (not numbers) in this workbook, so we rewrite the text in place:
#   - "1.234,56"  -> "1234.56"   (drop thousands dots, comma -> decimal point)
#   - "NAME A, NAME B" -> "NAME A. NAME B"   (stray comma -> period)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Razon social entries: comma used as a separator -> period ---
$ws.Range("E45").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E82").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"
$ws.Range("E46").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F46").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E47").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E83").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"

# --- "Importe" column (H2:H109): re-format the text values ---
# These are text cells holding numbers formatted like "1.005.000,00". A plain
# $ws.Range(...).Value = "1005000.00" assignment would be auto-coerced by Excel
# into the *number* 1005000, losing the intended text representation. Prefixing
# with a single quote forces Excel to store it as literal text (quote-prefix),
# exactly like typing an apostrophe before a numeric-looking entry. We then clear
# the quote-prefix formatting side effect so the cell keeps its original (default)
# style and only the shared-string text content changes.
$ws.Range("H2").Value = "'600.00"
$ws.Range("H3").Value = "'4832.50"
$ws.Range("H4").Value = "'700.00"
$ws.Range("H5").Value = "'400.00"
$ws.Range("H6").Value = "'798.60"
$ws.Range("H7").Value = "'278.00"
$ws.Range("H8").Value = "'77930.40"
$ws.Range("H9").Value = "'21869.79"
$ws.Range("H10").Value = "'5894.50"
$ws.Range("H11").Value = "'13196.55"
$ws.Range("H12").Value = "'2366.61"
$ws.Range("H13").Value = "'15316.11"
$ws.Range("H14").Value = "'9126.75"
$ws.Range("H15").Value = "'626.50"
$ws.Range("H16").Value = "'11273.11"
$ws.Range("H17").Value = "'900.00"
$ws.Range("H18").Value = "'26440.00"
$ws.Range("H19").Value = "'11641.82"
$ws.Range("H20").Value = "'391.50"
$ws.Range("H21").Value = "'6868.51"
$ws.Range("H22").Value = "'142.26"
$ws.Range("H23").Value = "'69080.64"
$ws.Range("H24").Value = "'286.03"
$ws.Range("H25").Value = "'18.00"
$ws.Range("H26").Value = "'1870.40"
$ws.Range("H27").Value = "'26218.26"
$ws.Range("H28").Value = "'18871.48"
$ws.Range("H29").Value = "'634.00"
$ws.Range("H30").Value = "'119.12"
$ws.Range("H31").Value = "'540.43"
$ws.Range("H32").Value = "'423.05"
$ws.Range("H33").Value = "'5299.00"
$ws.Range("H34").Value = "'4122.00"
$ws.Range("H35").Value = "'1680.00"
$ws.Range("H36").Value = "'14560.00"
$ws.Range("H37").Value = "'1504.51"
$ws.Range("H38").Value = "'37.09"
$ws.Range("H39").Value = "'2292.95"
$ws.Range("H40").Value = "'6610.00"
$ws.Range("H41").Value = "'1200.00"
$ws.Range("H42").Value = "'2444.00"
$ws.Range("H43").Value = "'18996.00"
$ws.Range("H44").Value = "'8130.00"
$ws.Range("H45").Value = "'1142.00"
$ws.Range("H46").Value = "'241.97"
$ws.Range("H47").Value = "'2070.00"
$ws.Range("H48").Value = "'190.00"
$ws.Range("H49").Value = "'0.04"
$ws.Range("H50").Value = "'157316.76"
$ws.Range("H51").Value = "'434610.00"
$ws.Range("H52").Value = "'0.92"
$ws.Range("H53").Value = "'3367.00"
$ws.Range("H54").Value = "'31.81"
$ws.Range("H55").Value = "'230.00"
$ws.Range("H56").Value = "'244.00"
$ws.Range("H57").Value = "'3860.00"
$ws.Range("H58").Value = "'1363.70"
$ws.Range("H59").Value = "'1339.00"
$ws.Range("H60").Value = "'810.00"
$ws.Range("H61").Value = "'16035.00"
$ws.Range("H62").Value = "'2398.00"
$ws.Range("H63").Value = "'788.45"
$ws.Range("H64").Value = "'29.56"
$ws.Range("H65").Value = "'69.00"
$ws.Range("H66").Value = "'1005000.00"
$ws.Range("H67").Value = "'2200.00"
$ws.Range("H68").Value = "'700.00"
$ws.Range("H69").Value = "'1000.00"
$ws.Range("H70").Value = "'2200.00"
$ws.Range("H71").Value = "'800.00"
$ws.Range("H72").Value = "'7823.37"
$ws.Range("H73").Value = "'500.00"
$ws.Range("H74").Value = "'950.00"
$ws.Range("H75").Value = "'650.00"
$ws.Range("H76").Value = "'1000.00"
$ws.Range("H77").Value = "'1600.00"
$ws.Range("H78").Value = "'200.00"
$ws.Range("H79").Value = "'240.00"
$ws.Range("H80").Value = "'2000.00"
$ws.Range("H81").Value = "'1000.00"
$ws.Range("H82").Value = "'250.00"
$ws.Range("H83").Value = "'120.00"
$ws.Range("H84").Value = "'260.00"
$ws.Range("H85").Value = "'105.00"
$ws.Range("H86").Value = "'1905.00"
$ws.Range("H87").Value = "'1470.97"
$ws.Range("H88").Value = "'2676.00"
$ws.Range("H89").Value = "'361.62"
$ws.Range("H90").Value = "'492.00"
$ws.Range("H91").Value = "'180.00"
$ws.Range("H92").Value = "'24.50"
$ws.Range("H93").Value = "'979.70"
$ws.Range("H94").Value = "'29960.00"
$ws.Range("H95").Value = "'500.00"
$ws.Range("H96").Value = "'877.20"
$ws.Range("H97").Value = "'119.76"
$ws.Range("H98").Value = "'2338.95"
$ws.Range("H99").Value = "'1760.00"
$ws.Range("H100").Value = "'1260.00"
$ws.Range("H101").Value = "'88500.00"
$ws.Range("H102").Value = "'65500.00"
$ws.Range("H103").Value = "'43000.00"
$ws.Range("H104").Value = "'47500.00"
$ws.Range("H105").Value = "'104000.00"
$ws.Range("H106").Value = "'86000.00"
$ws.Range("H107").Value = "'27500.00"
$ws.Range("H108").Value = "'417961.66"
$ws.Range("H109").Value = "'32400.00"

# Strip the quote-prefix style side effect back to the sheet default.
$ws.Range("H2:H109").ClearFormats()

